# Auto-generated edit script applying numeric corrections from the commit diff.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 22222392
$ws.Range("I11").Value = 22222392
$ws.Range("K11").Value = 22222392
$ws.Range("M11").Value = -22222252
$ws.Range("H51").Value = 3498.75
$ws.Range("J51").Value = 3998.3333
$ws.Range("L51").Value = 3998.3333
$ws.Range("N51").Value = -4966.3333
$ws.Range("H70").Value = 32975.25
$ws.Range("I70").Value = 950.5
$ws.Range("K70").Value = 2851.5
$ws.Range("M70").Value = -2581.5
$ws.Range("H73").Value = 32975.25
$ws.Range("I73").Value = 950.5
$ws.Range("K73").Value = 2851.5
$ws.Range("M73").Value = -1915.5
$ws.Range("H112").Value = 1486.6
$ws.Range("I112").Value = 999
$ws.Range("K112").Value = 2997
$ws.Range("M112").Value = -1889
$ws.Range("H116").Value = 12370.182
$ws.Range("J116").Value = 3635.25
$ws.Range("L116").Value = 3635.25
$ws.Range("N116").Value = -10519.25

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3534.1428
$ws.Range("I32").Value = 2789.7031
$ws.Range("J32").Value = 11474.833
$ws.Range("K32").Value = 2789.7031
$ws.Range("L32").Value = 11474.833
$ws.Range("M32").Value = -2502.7031
$ws.Range("N32").Value = -12048.833
$ws.Range("H53").Value = 8037
$ws.Range("I53").Value = 1074
$ws.Range("K53").Value = 1074
$ws.Range("M53").Value = -392
$ws.Range("H74").Value = 1319.25
$ws.Range("I74").Value = 478.78946
$ws.Range("K74").Value = 478.78946
$ws.Range("M74").Value = 395.21054
$ws.Range("H77").Value = 1319.25
$ws.Range("I77").Value = 478.78946
$ws.Range("K77").Value = 2393.9473
$ws.Range("M77").Value = 1974.0527
$ws.Range("H132").Value = 2431.9443
$ws.Range("I132").Value = 2051.8667
$ws.Range("K132").Value = 6155.6001
$ws.Range("M132").Value = -3625.6001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H22").Value = 1133.3334
$ws.Range("I22").Value = 460.2
$ws.Range("J22").Value = 1469.9
$ws.Range("K22").Value = 460.2
$ws.Range("L22").Value = 1469.9
$ws.Range("M22").Value = -110.2
$ws.Range("N22").Value = -2169.9
$ws.Range("H31").Value = 2181.4666
$ws.Range("I31").Value = 1736.8572
$ws.Range("J31").Value = 2570.5
$ws.Range("K31").Value = 1736.8572
$ws.Range("L31").Value = 2570.5
$ws.Range("M31").Value = -1441.8572
$ws.Range("N31").Value = -3160.5
$ws.Range("H34").Value = 2181.4666
$ws.Range("I34").Value = 1736.8572
$ws.Range("J34").Value = 2570.5
$ws.Range("K34").Value = 1736.8572
$ws.Range("L34").Value = 2570.5
$ws.Range("M34").Value = -1534.8572
$ws.Range("N34").Value = -2974.5
$ws.Range("H62").Value = 3650
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876
$ws.Range("H65").Value = 3650
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380
$ws.Range("H132").Value = 2146.2083
$ws.Range("I132").Value = 1214.8
$ws.Range("J132").Value = 3698.5557
$ws.Range("K132").Value = 3644.4
$ws.Range("L132").Value = 11095.6671
$ws.Range("M132").Value = -1114.4
$ws.Range("N132").Value = -16155.6671

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 498.2857
$ws.Range("I5").Value = 467.3846
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 1402.1538
$ws.Range("L5").Value = 2700
$ws.Range("M5").Value = -1290.1538
$ws.Range("N5").Value = -2924
$ws.Range("H75").Value = 985.2
$ws.Range("J75").Value = 985.2
$ws.Range("L75").Value = 2955.6
$ws.Range("N75").Value = -4951.6
$ws.Range("H78").Value = 985.2
$ws.Range("J78").Value = 985.2
$ws.Range("L78").Value = 8866.800000000001
$ws.Range("N78").Value = -18850.8
$ws.Range("H93").Value = 5999
$ws.Range("J93").Value = 5999
$ws.Range("L93").Value = 17997
$ws.Range("N93").Value = -21741
$ws.Range("H103").Value = 3169.5715
$ws.Range("I103").Value = 2512.5
$ws.Range("J103").Value = 3432.4
$ws.Range("K103").Value = 7537.5
$ws.Range("L103").Value = 10297.2
$ws.Range("M103").Value = -6658.5
$ws.Range("N103").Value = -12055.2
$ws.Range("H122").Value = 809.1111
$ws.Range("I122").Value = 642.6
$ws.Range("J122").Value = 1017.25
$ws.Range("K122").Value = 5783.400000000001
$ws.Range("L122").Value = 9155.25
$ws.Range("M122").Value = -3333.400000000001
$ws.Range("N122").Value = -14055.25
$ws.Range("H131").Value = 20198.75
$ws.Range("J131").Value = 23349.227
$ws.Range("L131").Value = 70047.681
$ws.Range("N131").Value = -80127.681
$ws.Range("H132").Value = 1601.8572
$ws.Range("I132").Value = 1042.6
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 9383.4
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -6853.4
$ws.Range("N132").Value = -32060
$ws.Range("H135").Value = 498.2857
$ws.Range("I135").Value = 467.3846
$ws.Range("J135").Value = 900
$ws.Range("K135").Value = 4206.4614
$ws.Range("L135").Value = 8100
$ws.Range("M135").Value = -1671.4614
$ws.Range("N135").Value = -13170

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 62750
$ws.Range("J5").Value = 62750
$ws.Range("L5").Value = 62750
$ws.Range("N5").Value = -62974
$ws.Range("H70").Value = 12766.667
$ws.Range("I70").Value = 23175
$ws.Range("K70").Value = 23175
$ws.Range("M70").Value = -22905
$ws.Range("H73").Value = 12766.667
$ws.Range("I73").Value = 23175
$ws.Range("K73").Value = 23175
$ws.Range("M73").Value = -22239
$ws.Range("H122").Value = 953.5333000000001
$ws.Range("I122").Value = 809.5
$ws.Range("J122").Value = 1118.1428
$ws.Range("K122").Value = 2428.5
$ws.Range("L122").Value = 3354.4284
$ws.Range("M122").Value = 21.5
$ws.Range("N122").Value = -8254.428400000001
$ws.Range("H132").Value = 1375678.4
$ws.Range("I132").Value = 1540099.8
$ws.Range("K132").Value = 4620299.4
$ws.Range("M132").Value = -4617769.4

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1825.2222
$ws.Range("I22").Value = 1289.4
$ws.Range("J22").Value = 2495
$ws.Range("K22").Value = 1289.4
$ws.Range("L22").Value = 2495
$ws.Range("M22").Value = -994.4000000000001
$ws.Range("N22").Value = -3085
$ws.Range("H27").Value = 1825.2222
$ws.Range("I27").Value = 1289.4
$ws.Range("J27").Value = 2495
$ws.Range("K27").Value = 1289.4
$ws.Range("L27").Value = 2495
$ws.Range("M27").Value = -1182.4
$ws.Range("N27").Value = -2709
$ws.Range("H46").Value = 1601.5294
$ws.Range("J46").Value = 2050.111
$ws.Range("L46").Value = 2050.111
$ws.Range("N46").Value = -2426.111
$ws.Range("H55").Value = 706.44446
$ws.Range("I55").Value = 899.6667
$ws.Range("J55").Value = 609.8333
$ws.Range("K55").Value = 899.6667
$ws.Range("L55").Value = 609.8333
$ws.Range("M55").Value = -726.6667
$ws.Range("N55").Value = -955.8333
$ws.Range("H100").Value = 1064.4445
$ws.Range("I100").Value = 1080.6666
$ws.Range("K100").Value = 1080.6666
$ws.Range("M100").Value = -539.6666
$ws.Range("H132").Value = 3112.6956
$ws.Range("I132").Value = 2112.375
$ws.Range("J132").Value = 3646.2
$ws.Range("K132").Value = 6337.125
$ws.Range("L132").Value = 10938.6
$ws.Range("M132").Value = -3807.125
$ws.Range("N132").Value = -15998.6

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H62").Value = 4398.4
$ws.Range("I62").Value = 3998
$ws.Range("K62").Value = 3998
$ws.Range("M62").Value = -3374
$ws.Range("H65").Value = 4398.4
$ws.Range("I65").Value = 3998
$ws.Range("K65").Value = 19990
$ws.Range("M65").Value = -16870
$ws.Range("H69").Value = 32499.5
$ws.Range("J69").Value = 32499.5
$ws.Range("L69").Value = 32499.5
$ws.Range("N69").Value = -33997.5
$ws.Range("H72").Value = 32499.5
$ws.Range("J72").Value = 32499.5
$ws.Range("L72").Value = 97498.5
$ws.Range("N72").Value = -104986.5
$ws.Range("H132").Value = 2278
$ws.Range("I132").Value = 1506.381
$ws.Range("K132").Value = 4519.143
$ws.Range("M132").Value = -1989.143
$ws.Range("H136").Value = 1353.3
$ws.Range("I136").Value = 939.05
$ws.Range("J136").Value = 2181.8
$ws.Range("K136").Value = 2817.15
$ws.Range("L136").Value = 6545.400000000001
$ws.Range("M136").Value = -267.1499999999996
$ws.Range("N136").Value = -11645.4
